$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Feature: Added Sub-heading for git setup and collaboration.
#
# Append a new bold sub-heading paragraph ("How to set up Git and work
# collaboratively with more than one person") after the last paragraph
# of the document ("A master branch in Git is the main branch...."),
# followed by a trailing blank paragraph (also carrying bold run
# formatting on its paragraph mark, matching the other section-heading
# blocks already present in the document, e.g. "Repository:", "Clone:",
# "Push", "Master Branch", etc.).
# ------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range
$tail.Collapse(0)                 # wdCollapseEnd
$tail.InsertParagraphAfter() | Out-Null

$headingIndex = $d.Paragraphs.Count
$heading = $d.Paragraphs.Item($headingIndex)

# Type the heading text as two runs ("H" then the remainder) to mirror
# the way the heading was authored (first character typed, then the
# rest of the phrase).
$headRange = $heading.Range
$headRange.Collapse(0)
$headRange.InsertAfter("H")
$headRange.Font.Bold = $true
$headRange.Font.BoldBi = $true
$headRange.Font.Size = 12

$restRange = $d.Paragraphs.Item($headingIndex).Range
$restRange.Collapse(0)
$restRange.MoveEnd(1, -1) | Out-Null   # wdCharacter, land before the pilcrow
$restRange.InsertAfter("ow to set up Git and work collaboratively with more than one person")
$restRange.Font.Bold = $true
$restRange.Font.BoldBi = $true
$restRange.Font.Size = 12

# Append the trailing empty paragraph after the new heading.
$afterHeading = $d.Paragraphs.Item($headingIndex).Range
$afterHeading.Collapse(0)
$afterHeading.InsertParagraphAfter() | Out-Null

$trailingIndex = $d.Paragraphs.Count
$trailing = $d.Paragraphs.Item($trailingIndex).Range
$trailing.Font.Bold = $true
$trailing.Font.BoldBi = $true
$trailing.Font.Size = 12

Write-Output ("Appended sub-heading paragraph " + $headingIndex + " and trailing paragraph " + $trailingIndex + ". Total paragraphs: " + $d.Paragraphs.Count)
